# Update cryptos list values (price + volume/1h) for Tue Feb 6 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most rows only change D (Price) and/or E (Volume 1h); rows 48/49 also swap
# their Coin name (B) and Link (C). D/E columns hold text-formatted values
# (e.g. "1.00", "301.78") that Excel would otherwise auto-convert to numbers,
# so force the cell format to Text ("@") before assigning.

function Set-Text($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-Text $ws.Range("D2") "42.999.03"
$ws.Range("E2").Value  = "  -0.08%  "

Set-Text $ws.Range("D3") "2.328.85"
$ws.Range("E3").Value  = "  +1.03%  "

Set-Text $ws.Range("D4") "1.00"
$ws.Range("E4").Value  = "  +0.15%  "

Set-Text $ws.Range("D5") "301.78"
$ws.Range("E5").Value  = "  -1.06%  "

Set-Text $ws.Range("D6") "95.98"
$ws.Range("E6").Value  = "  -1.20%  "

Set-Text $ws.Range("D7") "0.504"
$ws.Range("E7").Value  = "  -0.32%  "

$ws.Range("E8").Value  = "  +0.17%  "

$ws.Range("E9").Value  = "  -1.30%  "

Set-Text $ws.Range("D10") "34.39"
$ws.Range("E10").Value = "  -3.06%  "

Set-Text $ws.Range("D11") "19.08"
$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("E13").Value = "  +1.44%  "

Set-Text $ws.Range("D14") "6.76"
$ws.Range("E14").Value = "  -2.08%  "

Set-Text $ws.Range("D15") "2.695.32"
$ws.Range("E15").Value = "  +1.17%  "

Set-Text $ws.Range("D16") "2.338.35"
$ws.Range("E16").Value = "  +2.39%  "

Set-Text $ws.Range("D17") "0.792"
$ws.Range("E17").Value = "  +1.13%  "

Set-Text $ws.Range("D18") "42.943.81"
$ws.Range("E18").Value = "  +0.02%  "

Set-Text $ws.Range("D19") "12.28"
$ws.Range("E19").Value = "  -2.58%  "

Set-Text $ws.Range("D20") "6.18"
$ws.Range("E20").Value = "  +2.23%  "

Set-Text $ws.Range("D21") "0.0₃0893"
$ws.Range("E21").Value = "  -0.64%  "

Set-Text $ws.Range("D22") "68.01"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("E23").Value = "  +4.76%  "

Set-Text $ws.Range("D24") "236.71"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  -0.03%  "

Set-Text $ws.Range("D26") "2.40"
$ws.Range("E26").Value = "  -0.57%  "

Set-Text $ws.Range("D27") "24.70"
$ws.Range("E27").Value = "  -1.37%  "

Set-Text $ws.Range("D28") "2.06"
$ws.Range("E28").Value = "  -5.93%  "

$ws.Range("E29").Value = "  +1.10%  "

Set-Text $ws.Range("D30") "32.26"
$ws.Range("E30").Value = "  -2.16%  "

Set-Text $ws.Range("D31") "144.64"
$ws.Range("E31").Value = "  -12.94%  "

$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  +0.41%  "

Set-Text $ws.Range("D34") "17.87"
$ws.Range("E34").Value = "  -1.74%  "

$ws.Range("E35").Value = "  +1.92%  "

Set-Text $ws.Range("D36") "4.43"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  +2.96%  "

Set-Text $ws.Range("D38") "2.29"
$ws.Range("E38").Value = "  -2.46%  "

Set-Text $ws.Range("D39") "0.100"
$ws.Range("E39").Value = "  -0.58%  "

Set-Text $ws.Range("D40") "2.74"
$ws.Range("E40").Value = "  -0.01%  "

Set-Text $ws.Range("D41") "22.04"
$ws.Range("E41").Value = "  +23.38%  "

$ws.Range("E42").Value = "  -0.62%  "

Set-Text $ws.Range("D43") "1.934.11"
$ws.Range("E43").Value = "  -3.13%  "

$ws.Range("E44").Value = "  -0.44%  "

Set-Text $ws.Range("D45") "10.14"
$ws.Range("E45").Value = "  -2.82%  "

$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("E47").Value = "  -1.18%  "

# Rows 48/49 swap coin identity (name + link); new price/volume values follow.
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-Text $ws.Range("D48") "2.561.64"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-Text $ws.Range("D49") "2.87"
$ws.Range("E49").Value = "  +1.12%  "

Set-Text $ws.Range("D50") "53.70"
$ws.Range("E50").Value = "  +0.18%  "

Set-Text $ws.Range("D51") "73.14"
$ws.Range("E51").Value = "  +2.03%  "
